# Insert a new data row above the existing row 31 ("Fecha" = 2023-06-02 / serial 45079),
# which pushes the former rows 31-49 down to rows 32-50, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(31).Insert()

$ws.Range("A31").Value = 9
$ws.Range("B31").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C31").Value = "Metropolitana"
$ws.Range("D31").Value = 45079
$ws.Range("D31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E31").Value = 13
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100104
$ws.Range("H31").Value = "Frutos de pepita"
$ws.Range("I31").Value = 100104001
$ws.Range("J31").Value = "Granada"
$ws.Range("K31").Value = "Wonderfull"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 200
$ws.Range("N31").Value = 7500
$ws.Range("O31").Value = 7500
$ws.Range("P31").Value = 7500
$ws.Range("Q31").Value = "`$/caja 15 kilos granel"
$ws.Range("R31").Value = "Provincia de Los Andes"
$ws.Range("S31").Value = 500
$ws.Range("T31").Value = 15
